# Apply "Ventas objetivo" formula change: update computed values in
# columns L (Diferencia Stock), R (uds. Objetivo semana pasada),
# T (Tendencia Consumo) and the Total_Ajuste_Stock summary cell C146.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_7")

$ws.Range("R11").Value = 2
$ws.Range("R13").Value = 1
$ws.Range("R14").Value = 3
$ws.Range("R24").Value = 2
$ws.Range("L25").Value = 0
$ws.Range("R25").Value = 3
$ws.Range("T25").Value = 1
$ws.Range("R37").Value = 1
$ws.Range("R38").Value = 1
$ws.Range("R45").Value = 1
$ws.Range("R54").Value = 2
$ws.Range("R55").Value = 1
$ws.Range("L56").Value = 0
$ws.Range("R56").Value = 2
$ws.Range("R73").Value = 1
$ws.Range("R74").Value = 1
$ws.Range("L75").Value = 0
$ws.Range("R75").Value = 4
$ws.Range("R79").Value = 2
$ws.Range("T79").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("R85").Value = 6
$ws.Range("L95").Value = 0
$ws.Range("R95").Value = 1
$ws.Range("L101").Value = 0
$ws.Range("R101").Value = 10
$ws.Range("R103").Value = 2
$ws.Range("R105").Value = 1
$ws.Range("R111").Value = 1
$ws.Range("T111").Value = 1
$ws.Range("L114").Value = 0
$ws.Range("R114").Value = 4
$ws.Range("R115").Value = 1
$ws.Range("R122").Value = 1
$ws.Range("R124").Value = 1
$ws.Range("T124").Value = 0
$ws.Range("R130").Value = 1
$ws.Range("C146").Value = 0
